# Insert a new weekly price record for "Berenjena" at Vega Monumental
# Concepción. The new observation (Fecha = 45027) belongs between the
# existing row 106 (Fecha = 45016) and the old row 107 (Fecha = 44971),
# so a whole row is inserted at position 107, pushing every following
# row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 107 — shifts old rows 107..148 down to 108..149
# and carries formatting (so column D keeps its date/time style).
$ws.Rows.Item(107).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(107, 1).Value  = 11
$ws.Cells.Item(107, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(107, 3).Value  = "Bíobío"
$ws.Cells.Item(107, 4).Value  = 45027
$ws.Cells.Item(107, 5).Value  = 8
$ws.Cells.Item(107, 6).Value  = 100112001
$ws.Cells.Item(107, 7).Value  = "Berenjena"
$ws.Cells.Item(107, 8).Value  = "Sin especificar"
$ws.Cells.Item(107, 9).Value  = "Primera"
$ws.Cells.Item(107, 10).Value = 180
$ws.Cells.Item(107, 11).Value = 7000
$ws.Cells.Item(107, 12).Value = 7500
$ws.Cells.Item(107, 13).Value = 7278
$ws.Cells.Item(107, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(107, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(107, 16).Value = 121
$ws.Cells.Item(107, 17).Value = 60
$ws.Cells.Item(107, 18).Value = "Hortaliza"
